$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.684.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.630.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0783"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.629.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.854.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.711.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.544"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.104.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.795"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.762.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0109"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("B48").Value = "SynthetixNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.13%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.419"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.22%  "
